$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new "Region" column header and related info, mirroring the
# existing "Active"/Boolean column pattern but for an Enum "Region".
$ws.Range("H1").Value = "Region"
$ws.Range("H4").Value = "Enum"
$ws.Range("H5").Value = "North, South, etc."

# Set the new column's width to match the rest of the table style
# (closest achievable value to the source column width of 18.140625).
$ws.Columns.Item(8).ColumnWidth = 17.3

# Update the active selection to the newly added cell, as in the diff.
$ws.Range("H5").Select()
